$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column E/F: hide them
# ---------------------------------------------------------------------------
$ws.Columns("E:F").Hidden = $true

# ---------------------------------------------------------------------------
# Helper dates (explicit midnight to avoid fractional serials)
# ---------------------------------------------------------------------------
$d0920 = Get-Date -Year 2023 -Month 9 -Day 20 -Hour 0 -Minute 0 -Second 0
$d0923 = Get-Date -Year 2023 -Month 9 -Day 23 -Hour 0 -Minute 0 -Second 0
$d0925 = Get-Date -Year 2023 -Month 9 -Day 25 -Hour 0 -Minute 0 -Second 0
$d0926 = Get-Date -Year 2023 -Month 9 -Day 26 -Hour 0 -Minute 0 -Second 0

# ---------------------------------------------------------------------------
# Row 7: update date + mark "Si " (trailing space) in G7
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = $d0920
$ws.Range("G7").Value = "Si "

# ---------------------------------------------------------------------------
# Row 8: new description text, update date, mark "Si" in G8
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "Registrar usuarios, loguin y logout"
$ws.Range("D8").Value = $d0923
$ws.Range("G8").Value = "Si"
$ws.Rows(8).RowHeight = 57.6

# ---------------------------------------------------------------------------
# Row 9: new description text (B9 previously empty), update date, mark G9
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = "Revisar que sin loguearse no acceda a carga de libros, modicacion ni eliminacino"
$ws.Range("D9").Value = $d0923
$ws.Range("G9").Value = "Si"
$ws.Rows(9).RowHeight = 57.6

# Give B9 the same left/center/wrap formatting style as B7/B8 (it was a plain
# empty cell before; now needs the "description" look).
$rngB9 = $ws.Range("B9")
$rngB9.HorizontalAlignment = -4131   # xlLeft
$rngB9.VerticalAlignment = -4108     # xlCenter
$rngB9.WrapText = $true
$rngB9.Borders.LineStyle = 1
$rngB9.Borders.Weight = 2
$rngB9.Borders.Color = 0

Write-Output "rows 7-9 updated"

# ---------------------------------------------------------------------------
# New rows 10-14: five additional test cases. Build each one with the same
# visual layout as the existing rows (center-aligned ID/Si columns, left
# aligned wrapped description, bordered A:G, border-less H:Z, date format).
# ---------------------------------------------------------------------------
function Format-TestCaseRow([int]$r, [double]$height) {
    $ws.Rows($r).RowHeight = $height

    $rngA = $ws.Range("A$r")
    $rngA.HorizontalAlignment = -4108   # xlCenter
    $rngA.VerticalAlignment = -4108     # xlCenter
    $rngA.WrapText = $false
    $rngA.Borders.LineStyle = 1
    $rngA.Borders.Weight = 2
    $rngA.Borders.Color = 0

    $rngG = $ws.Range("G$r")
    $rngG.HorizontalAlignment = -4108   # xlCenter
    $rngG.VerticalAlignment = -4108     # xlCenter
    $rngG.WrapText = $false
    $rngG.Borders.LineStyle = 1
    $rngG.Borders.Weight = 2
    $rngG.Borders.Color = 0

    $rngB = $ws.Range("B$r")
    $rngB.HorizontalAlignment = -4131   # xlLeft
    $rngB.VerticalAlignment = -4108     # xlCenter
    $rngB.WrapText = $true
    $rngB.Borders.LineStyle = 1
    $rngB.Borders.Weight = 2
    $rngB.Borders.Color = 0

    $rngC = $ws.Range("C$r")
    $rngC.HorizontalAlignment = -4131   # xlLeft
    $rngC.VerticalAlignment = -4108     # xlCenter
    $rngC.WrapText = $true
    $rngC.Borders.LineStyle = 1
    $rngC.Borders.Weight = 2
    $rngC.Borders.Color = 0

    $rngD = $ws.Range("D$r")
    $rngD.HorizontalAlignment = -4108   # xlCenter
    $rngD.VerticalAlignment = -4108     # xlCenter
    $rngD.WrapText = $true
    $rngD.Borders.LineStyle = 1
    $rngD.Borders.Weight = 2
    $rngD.Borders.Color = 0
    $rngD.NumberFormat = "d-mmm"

    $rngEF = $ws.Range("E$r:F$r")
    $rngEF.HorizontalAlignment = -4108  # xlCenter
    $rngEF.VerticalAlignment = -4108    # xlCenter
    $rngEF.WrapText = $true
    $rngEF.Borders.LineStyle = 1
    $rngEF.Borders.Weight = 2
    $rngEF.Borders.Color = 0

    $rngRest = $ws.Range("H$r:Z$r")
    $rngRest.HorizontalAlignment = -4108  # xlCenter
    $rngRest.VerticalAlignment = -4108    # xlCenter
    $rngRest.WrapText = $false
    $rngRest.Borders.LineStyle = -4142    # xlLineStyleNone
}

Format-TestCaseRow 10 54
Format-TestCaseRow 11 54
Format-TestCaseRow 12 54
Format-TestCaseRow 13 54
Format-TestCaseRow 14 54

$ws.Range("A10").Value = "Caso #4"
$ws.Range("B10").Value = "Nuevo libro, revisar que aparerza en el listado. Modificar y luego verificar que que aparezcan los nuevos datos. Pruebo eliminar un libro"
$ws.Range("D10").Value = $d0923
$ws.Range("G10").Value = "Si"

$ws.Range("A11").Value = "Caso #5"
$ws.Range("B11").Value = "Reservar libro. Revisar que se envie un mail al lugar con copia al usuario. "
$ws.Range("D11").Value = $d0925
$ws.Range("G11").Value = "Si"

$ws.Range("A12").Value = "Caso #6"
$ws.Range("B12").Value = "Listado de libros disponibles, que no aparezcan los reservados"
$ws.Range("D12").Value = $d0925
$ws.Range("G12").Value = "Si"

$ws.Range("A13").Value = "Caso #7"
$ws.Range("B13").Value = "Cargar comentario sobre un libro y que aparezca en el listado debajo del libro"
$ws.Range("D13").Value = $d0925
$ws.Range("G13").Value = "Si"

$ws.Range("A14").Value = "Caso #8"
$ws.Range("B14").Value = "Listado de libros reservados"
$ws.Range("D14").Value = $d0926
$ws.Range("G14").Value = "Si"

Write-Output "rows 10-14 added"
